$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 4
$ws.Range("I2").Value = "target"
$ws.Range("J2").Value = "old"
$ws.Range("K2").Value = "j"
$ws.Range("L2").Value = "stimuli/img_gqy6z.png"
$ws.Range("M2").Value = 86.47368421052632
$ws.Range("N2").Value = 68.42105263157895
$ws.Range("O2").Value = 77.44736842105263
$ws.Range("P2").Value = 38
$ws.Range("Q2").Value = 9
$ws.Range("R2").Value = 9
$ws.Range("S2").Value = 9
$ws.Range("T2").Value = 9
$ws.Range("U2").Value = 9
$ws.Range("V2").Value = 9

# Row 3
$ws.Range("C3").Value = 4
$ws.Range("L3").Value = "stimuli/img_zt893.png"
$ws.Range("M3").Value = 68.53191489361703
$ws.Range("N3").Value = 49.19148936170212
$ws.Range("O3").Value = 58.86170212765958
$ws.Range("P3").Value = 47
$ws.Range("Q3").Value = 5
$ws.Range("R3").Value = 5
$ws.Range("S3").Value = 5
$ws.Range("T3").Value = 5
$ws.Range("U3").Value = 5
$ws.Range("V3").Value = 5

# Row 4
$ws.Range("C4").Value = 4
$ws.Range("L4").Value = "stimuli/img_a8y4y.png"
$ws.Range("M4").Value = 75.15789473684211
$ws.Range("N4").Value = 53.76315789473684
$ws.Range("O4").Value = 64.46052631578948
$ws.Range("P4").Value = 38
$ws.Range("Q4").Value = 6
$ws.Range("R4").Value = 6
$ws.Range("S4").Value = 6
$ws.Range("T4").Value = 6
$ws.Range("U4").Value = 6
$ws.Range("V4").Value = 6

# Row 5
$ws.Range("C5").Value = 4
$ws.Range("I5").Value = "target"
$ws.Range("J5").Value = "old"
$ws.Range("K5").Value = "j"
$ws.Range("L5").Value = "stimuli/img_5yhyk.png"
$ws.Range("M5").Value = 46.375
$ws.Range("N5").Value = 31.325
$ws.Range("O5").Value = 38.85
$ws.Range("P5").Value = 40
$ws.Range("Q5").Value = 2
$ws.Range("R5").Value = 2
$ws.Range("S5").Value = 2
$ws.Range("T5").Value = 2
$ws.Range("U5").Value = 2
$ws.Range("V5").Value = 2

# Row 6
$ws.Range("C6").Value = 4
$ws.Range("I6").Value = $null
$ws.Range("J6").Value = "new"
$ws.Range("K6").Value = "f"
$ws.Range("L6").Value = "stimuli/img_2js6m.png"
$ws.Range("M6").Value = 40.02777777777778
$ws.Range("N6").Value = 20.88888888888889
$ws.Range("O6").Value = 30.45833333333334
$ws.Range("P6").Value = 36
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 2
$ws.Range("S6").Value = 2
$ws.Range("T6").Value = 2
$ws.Range("U6").Value = 2
$ws.Range("V6").Value = 2

# Row 7
$ws.Range("C7").Value = 4
$ws.Range("L7").Value = "stimuli/img_4wq98.png"
$ws.Range("M7").Value = 78.48387096774194
$ws.Range("N7").Value = 58.12903225806452
$ws.Range("O7").Value = 68.30645161290323
$ws.Range("P7").Value = 31

# Row 8
$ws.Range("C8").Value = 4
$ws.Range("I8").Value = "target"
$ws.Range("J8").Value = "old"
$ws.Range("K8").Value = "j"
$ws.Range("L8").Value = "stimuli/img_0eflx.png"
$ws.Range("M8").Value = 76.05128205128206
$ws.Range("N8").Value = 53.53846153846154
$ws.Range("O8").Value = 64.7948717948718
$ws.Range("P8").Value = 39
$ws.Range("Q8").Value = 6
$ws.Range("R8").Value = 6
$ws.Range("S8").Value = 6
$ws.Range("T8").Value = 6
$ws.Range("U8").Value = 6
$ws.Range("V8").Value = 6

# Row 9
$ws.Range("C9").Value = 4
$ws.Range("L9").Value = "stimuli/img_0nckg.png"
$ws.Range("M9").Value = 65.94285714285714
$ws.Range("N9").Value = 41.17142857142857
$ws.Range("O9").Value = 53.55714285714285
$ws.Range("P9").Value = 35
$ws.Range("Q9").Value = 4
$ws.Range("R9").Value = 4
$ws.Range("S9").Value = 4
$ws.Range("T9").Value = 4
$ws.Range("U9").Value = 4
$ws.Range("V9").Value = 4

# Row 10
$ws.Range("C10").Value = 4
$ws.Range("L10").Value = "stimuli/img_psgf7.png"
$ws.Range("M10").Value = 26
$ws.Range("N10").Value = 11.66666666666667
$ws.Range("O10").Value = 18.83333333333333
$ws.Range("P10").Value = 36
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = 1
$ws.Range("S10").Value = 1
$ws.Range("T10").Value = 1
$ws.Range("U10").Value = 1
$ws.Range("V10").Value = 1

# Row 11
$ws.Range("C11").Value = 4
$ws.Range("I11").Value = "target"
$ws.Range("J11").Value = "old"
$ws.Range("K11").Value = "j"
$ws.Range("L11").Value = "stimuli/img_cmyvx.png"
$ws.Range("M11").Value = 64.25
$ws.Range("N11").Value = 40.09375
$ws.Range("O11").Value = 52.171875
$ws.Range("P11").Value = 32
$ws.Range("Q11").Value = 4
$ws.Range("R11").Value = 4
$ws.Range("S11").Value = 4
$ws.Range("T11").Value = 4
$ws.Range("U11").Value = 4
$ws.Range("V11").Value = 4

# Row 12
$ws.Range("C12").Value = 4
$ws.Range("L12").Value = "stimuli/img_juob3.png"
$ws.Range("M12").Value = 79.92105263157895
$ws.Range("N12").Value = 59.78947368421053
$ws.Range("O12").Value = 69.85526315789474
$ws.Range("P12").Value = 38

# Row 13
$ws.Range("C13").Value = 4
$ws.Range("L13").Value = "stimuli/img_3h4c9.png"
$ws.Range("M13").Value = 85.47619047619048
$ws.Range("N13").Value = 67.26190476190476
$ws.Range("O13").Value = 76.36904761904762
$ws.Range("P13").Value = 42
$ws.Range("Q13").Value = 9
$ws.Range("R13").Value = 9
$ws.Range("S13").Value = 9
$ws.Range("T13").Value = 9
$ws.Range("U13").Value = 9
$ws.Range("V13").Value = 9

# Row 14
$ws.Range("C14").Value = 4
$ws.Range("I14").Value = "target"
$ws.Range("J14").Value = "old"
$ws.Range("K14").Value = "j"
$ws.Range("L14").Value = "stimuli/img_scrdm.png"
$ws.Range("M14").Value = 78.675
$ws.Range("N14").Value = 57.9
$ws.Range("O14").Value = 68.28749999999999
$ws.Range("P14").Value = 40
$ws.Range("Q14").Value = 7
$ws.Range("R14").Value = 7
$ws.Range("S14").Value = 7
$ws.Range("T14").Value = 7
$ws.Range("U14").Value = 7
$ws.Range("V14").Value = 7

# Row 15
$ws.Range("C15").Value = 4
$ws.Range("L15").Value = "stimuli/img_wyctg.png"
$ws.Range("M15").Value = 33.44736842105263
$ws.Range("N15").Value = 11.39473684210526
$ws.Range("O15").Value = 22.42105263157895
$ws.Range("P15").Value = 38
$ws.Range("Q15").Value = 1
$ws.Range("R15").Value = 1
$ws.Range("S15").Value = 1
$ws.Range("T15").Value = 1
$ws.Range("U15").Value = 1
$ws.Range("V15").Value = 1

# Row 16
$ws.Range("C16").Value = 4
$ws.Range("I16").Value = $null
$ws.Range("J16").Value = "new"
$ws.Range("K16").Value = "f"
$ws.Range("L16").Value = "stimuli/img_oou46.png"
$ws.Range("M16").Value = 75.70270270270271
$ws.Range("N16").Value = 54.86486486486486
$ws.Range("O16").Value = 65.28378378378379
$ws.Range("P16").Value = 37
$ws.Range("Q16").Value = 6
$ws.Range("R16").Value = 6
$ws.Range("S16").Value = 6
$ws.Range("T16").Value = 6
$ws.Range("U16").Value = 6
$ws.Range("V16").Value = 6

# Row 17
$ws.Range("C17").Value = 4
$ws.Range("I17").Value = "target"
$ws.Range("J17").Value = "old"
$ws.Range("K17").Value = "j"
$ws.Range("L17").Value = "stimuli/img_ca8kd.png"
$ws.Range("M17").Value = 92.05405405405405
$ws.Range("N17").Value = 73.02702702702703
$ws.Range("O17").Value = 82.54054054054055
$ws.Range("P17").Value = 37
$ws.Range("Q17").Value = 10
$ws.Range("R17").Value = 10
$ws.Range("S17").Value = 10
$ws.Range("T17").Value = 10
$ws.Range("U17").Value = 10
$ws.Range("V17").Value = 10

# Row 18
$ws.Range("C18").Value = 4
$ws.Range("L18").Value = "stimuli/img_rvssl.png"
$ws.Range("M18").Value = 74.25
$ws.Range("N18").Value = 54.33333333333334
$ws.Range("O18").Value = 64.29166666666667
$ws.Range("P18").Value = 36

# Row 19
$ws.Range("C19").Value = 4
$ws.Range("L19").Value = "stimuli/img_72fmj.png"
$ws.Range("M19").Value = 53.87179487179487
$ws.Range("N19").Value = 36.02564102564103
$ws.Range("O19").Value = 44.94871794871795
$ws.Range("P19").Value = 39
$ws.Range("Q19").Value = 3
$ws.Range("R19").Value = 3
$ws.Range("S19").Value = 3
$ws.Range("T19").Value = 3
$ws.Range("U19").Value = 3
$ws.Range("V19").Value = 3

# Row 20
$ws.Range("C20").Value = 4
$ws.Range("L20").Value = "stimuli/img_ozxpp.png"
$ws.Range("M20").Value = 26.26470588235294
$ws.Range("N20").Value = 11.47058823529412
$ws.Range("O20").Value = 18.86764705882353
$ws.Range("P20").Value = 34
$ws.Range("Q20").Value = 1
$ws.Range("R20").Value = 1
$ws.Range("S20").Value = 1
$ws.Range("T20").Value = 1
$ws.Range("U20").Value = 1
$ws.Range("V20").Value = 1

# Row 21
$ws.Range("C21").Value = 4
$ws.Range("L21").Value = "stimuli/img_ce55l.png"
$ws.Range("M21").Value = 82.23809523809524
$ws.Range("N21").Value = 64.07142857142857
$ws.Range("O21").Value = 73.1547619047619
$ws.Range("P21").Value = 42
$ws.Range("Q21").Value = 8
$ws.Range("R21").Value = 8
$ws.Range("S21").Value = 8
$ws.Range("T21").Value = 8
$ws.Range("U21").Value = 8
$ws.Range("V21").Value = 8

# Row 22
$ws.Range("C22").Value = 4
$ws.Range("I22").Value = "target"
$ws.Range("J22").Value = "old"
$ws.Range("K22").Value = "j"
$ws.Range("L22").Value = "stimuli/img_zi682.png"
$ws.Range("M22").Value = 84.59999999999999
$ws.Range("N22").Value = 69.52500000000001
$ws.Range("O22").Value = 77.0625
$ws.Range("P22").Value = 40
$ws.Range("Q22").Value = 9
$ws.Range("R22").Value = 9
$ws.Range("S22").Value = 9
$ws.Range("T22").Value = 9
$ws.Range("U22").Value = 9
$ws.Range("V22").Value = 9

# Row 23
$ws.Range("C23").Value = 4
$ws.Range("H23").Value = "bedrooms"
$ws.Range("J23").Value = "new"
$ws.Range("L23").Value = "stimuli/img_f4jxo.png"
$ws.Range("M23").Value = 82.91666666666667
$ws.Range("N23").Value = 65.52777777777777
$ws.Range("O23").Value = 74.22222222222223
$ws.Range("P23").Value = 36
$ws.Range("Q23").Value = 8
$ws.Range("R23").Value = 8
$ws.Range("S23").Value = 8
$ws.Range("T23").Value = 8
$ws.Range("U23").Value = 8
$ws.Range("V23").Value = 8

# Row 24
$ws.Range("C24").Value = 4
$ws.Range("I24").Value = "target"
$ws.Range("J24").Value = "old"
$ws.Range("K24").Value = "j"
$ws.Range("L24").Value = "stimuli/img_wijef.png"
$ws.Range("M24").Value = 69.875
$ws.Range("N24").Value = 48.025
$ws.Range("O24").Value = 58.95
$ws.Range("P24").Value = 40
$ws.Range("Q24").Value = 5
$ws.Range("R24").Value = 5
$ws.Range("S24").Value = 5
$ws.Range("T24").Value = 5
$ws.Range("U24").Value = 5
$ws.Range("V24").Value = 5

# Row 25
$ws.Range("C25").Value = 4
$ws.Range("I25").Value = $null
$ws.Range("J25").Value = "new"
$ws.Range("K25").Value = "f"
$ws.Range("L25").Value = "stimuli/img_1vq1v.png"
$ws.Range("M25").Value = 69.42857142857143
$ws.Range("N25").Value = 46.59523809523809
$ws.Range("O25").Value = 58.01190476190476
$ws.Range("Q25").Value = 5
$ws.Range("R25").Value = 5
$ws.Range("S25").Value = 5
$ws.Range("T25").Value = 5
$ws.Range("U25").Value = 5
$ws.Range("V25").Value = 5

# Row 26
$ws.Range("C26").Value = 4
$ws.Range("I26").Value = "target"
$ws.Range("J26").Value = "old"
$ws.Range("K26").Value = "j"
$ws.Range("L26").Value = "stimuli/img_c4uwt.png"
$ws.Range("M26").Value = 44.48387096774194
$ws.Range("N26").Value = 30.06451612903226
$ws.Range("O26").Value = 37.2741935483871
$ws.Range("P26").Value = 31
$ws.Range("Q26").Value = 2
$ws.Range("R26").Value = 2
$ws.Range("S26").Value = 2
$ws.Range("T26").Value = 2
$ws.Range("U26").Value = 2
$ws.Range("V26").Value = 2

# Row 27
$ws.Range("C27").Value = 4
$ws.Range("I27").Value = $null
$ws.Range("J27").Value = "new"
$ws.Range("K27").Value = "f"
$ws.Range("L27").Value = "stimuli/img_uxxo0.png"
$ws.Range("M27").Value = 71.74418604651163
$ws.Range("N27").Value = 48.44186046511628
$ws.Range("O27").Value = 60.09302325581395
$ws.Range("P27").Value = 43
$ws.Range("Q27").Value = 5
$ws.Range("R27").Value = 5
$ws.Range("S27").Value = 5
$ws.Range("T27").Value = 5
$ws.Range("U27").Value = 5
$ws.Range("V27").Value = 5

# Row 28
$ws.Range("C28").Value = 4
$ws.Range("L28").Value = "stimuli/img_t2ioc.png"
$ws.Range("M28").Value = 88.18918918918919
$ws.Range("N28").Value = 74.05405405405405
$ws.Range("O28").Value = 81.12162162162161
$ws.Range("P28").Value = 37
$ws.Range("Q28").Value = 10
$ws.Range("R28").Value = 10
$ws.Range("S28").Value = 10
$ws.Range("T28").Value = 10
$ws.Range("U28").Value = 10
$ws.Range("V28").Value = 10

# Row 29
$ws.Range("C29").Value = 4
$ws.Range("H29").Value = $null
$ws.Range("J29").Value = "catch"
$ws.Range("L29").Value = "stimuli/catch_24.jpg"
$ws.Range("M29").Value = $null
$ws.Range("N29").Value = $null
$ws.Range("O29").Value = $null
$ws.Range("P29").Value = $null
$ws.Range("Q29").Value = $null
$ws.Range("R29").Value = $null
$ws.Range("S29").Value = $null
$ws.Range("T29").Value = $null
$ws.Range("U29").Value = $null
$ws.Range("V29").Value = $null

# Row 30
$ws.Range("C30").Value = 4
$ws.Range("I30").Value = $null
$ws.Range("J30").Value = "new"
$ws.Range("K30").Value = "f"
$ws.Range("L30").Value = "stimuli/img_x0u5z.png"
$ws.Range("M30").Value = 92
$ws.Range("N30").Value = 78.16216216216216
$ws.Range("O30").Value = 85.08108108108108
$ws.Range("P30").Value = 37
$ws.Range("Q30").Value = 10
$ws.Range("R30").Value = 10
$ws.Range("S30").Value = 10
$ws.Range("T30").Value = 10
$ws.Range("U30").Value = 10
$ws.Range("V30").Value = 10

# Row 31
$ws.Range("C31").Value = 4
$ws.Range("I31").Value = $null
$ws.Range("J31").Value = "new"
$ws.Range("K31").Value = "f"
$ws.Range("L31").Value = "stimuli/img_a9acb.png"
$ws.Range("M31").Value = 77.11428571428571
$ws.Range("N31").Value = 58.42857142857143
$ws.Range("O31").Value = 67.77142857142857
$ws.Range("P31").Value = 35
$ws.Range("Q31").Value = 7
$ws.Range("R31").Value = 7
$ws.Range("S31").Value = 7
$ws.Range("T31").Value = 7
$ws.Range("U31").Value = 7
$ws.Range("V31").Value = 7

# Row 32
$ws.Range("C32").Value = 4
$ws.Range("I32").Value = "target"
$ws.Range("J32").Value = "old"
$ws.Range("K32").Value = "j"
$ws.Range("L32").Value = "stimuli/img_bpyv5.png"
$ws.Range("M32").Value = 59.05882352941177
$ws.Range("N32").Value = 37.55882352941177
$ws.Range("O32").Value = 48.30882352941177
$ws.Range("P32").Value = 34
$ws.Range("Q32").Value = 3
$ws.Range("R32").Value = 3
$ws.Range("S32").Value = 3
$ws.Range("T32").Value = 3
$ws.Range("U32").Value = 3
$ws.Range("V32").Value = 3

# Row 33
$ws.Range("C33").Value = 4
$ws.Range("I33").Value = "target"
$ws.Range("J33").Value = "old"
$ws.Range("K33").Value = "j"
$ws.Range("L33").Value = "stimuli/img_g2akb.png"
$ws.Range("M33").Value = 87.875
$ws.Range("N33").Value = 79
$ws.Range("O33").Value = 83.4375
$ws.Range("Q33").Value = 10
$ws.Range("R33").Value = 10
$ws.Range("S33").Value = 10
$ws.Range("T33").Value = 10
$ws.Range("U33").Value = 10
$ws.Range("V33").Value = 10

# Row 34
$ws.Range("C34").Value = 4
$ws.Range("I34").Value = "target"
$ws.Range("J34").Value = "old"
$ws.Range("K34").Value = "j"
$ws.Range("L34").Value = "stimuli/img_kljj4.png"
$ws.Range("M34").Value = 64.34999999999999
$ws.Range("N34").Value = 44.15
$ws.Range("O34").Value = 54.25
$ws.Range("P34").Value = 40
$ws.Range("Q34").Value = 4
$ws.Range("R34").Value = 4
$ws.Range("S34").Value = 4
$ws.Range("T34").Value = 4
$ws.Range("U34").Value = 4
$ws.Range("V34").Value = 4

# Row 35
$ws.Range("C35").Value = 4
$ws.Range("I35").Value = $null
$ws.Range("J35").Value = "new"
$ws.Range("K35").Value = "f"
$ws.Range("L35").Value = "stimuli/img_zgg62.png"
$ws.Range("M35").Value = 82.18421052631579
$ws.Range("N35").Value = 63.52631578947368
$ws.Range("O35").Value = 72.85526315789474
$ws.Range("P35").Value = 38
$ws.Range("Q35").Value = 8
$ws.Range("R35").Value = 8
$ws.Range("S35").Value = 8
$ws.Range("T35").Value = 8
$ws.Range("U35").Value = 8
$ws.Range("V35").Value = 8

# Row 36
$ws.Range("C36").Value = 4
$ws.Range("I36").Value = $null
$ws.Range("J36").Value = "new"
$ws.Range("K36").Value = "f"
$ws.Range("L36").Value = "stimuli/img_u2o6z.png"
$ws.Range("M36").Value = 58.6
$ws.Range("N36").Value = 38.2
$ws.Range("O36").Value = 48.40000000000001
$ws.Range("P36").Value = 30
$ws.Range("Q36").Value = 3
$ws.Range("R36").Value = 3
$ws.Range("S36").Value = 3
$ws.Range("T36").Value = 3
$ws.Range("U36").Value = 3
$ws.Range("V36").Value = 3

# Row 37
$ws.Range("C37").Value = 4
$ws.Range("I37").Value = $null
$ws.Range("J37").Value = "new"
$ws.Range("K37").Value = "f"
$ws.Range("L37").Value = "stimuli/img_jp28n.png"
$ws.Range("M37").Value = 65.02564102564102
$ws.Range("N37").Value = 44.97435897435897
$ws.Range("O37").Value = 55
$ws.Range("P37").Value = 39
$ws.Range("Q37").Value = 4
$ws.Range("R37").Value = 4
$ws.Range("S37").Value = 4
$ws.Range("T37").Value = 4
$ws.Range("U37").Value = 4
$ws.Range("V37").Value = 5

# Row 38
$ws.Range("C38").Value = 4
$ws.Range("I38").Value = $null
$ws.Range("J38").Value = "new"
$ws.Range("K38").Value = "f"
$ws.Range("L38").Value = "stimuli/img_le8uf.png"
$ws.Range("M38").Value = 12.88888888888889
$ws.Range("N38").Value = 9.222222222222221
$ws.Range("O38").Value = 11.05555555555556
$ws.Range("P38").Value = 36
$ws.Range("Q38").Value = 1
$ws.Range("R38").Value = 1
$ws.Range("S38").Value = 1
$ws.Range("T38").Value = 1
$ws.Range("U38").Value = 1
$ws.Range("V38").Value = 1

# Row 39
$ws.Range("C39").Value = 4
$ws.Range("I39").Value = $null
$ws.Range("J39").Value = "new"
$ws.Range("K39").Value = "f"
$ws.Range("L39").Value = "stimuli/img_e26ut.png"
$ws.Range("M39").Value = 81.07692307692308
$ws.Range("N39").Value = 61.28205128205128
$ws.Range("O39").Value = 71.17948717948718
$ws.Range("P39").Value = 39
$ws.Range("Q39").Value = 8
$ws.Range("R39").Value = 8
$ws.Range("S39").Value = 8
$ws.Range("T39").Value = 8
$ws.Range("U39").Value = 8
$ws.Range("V39").Value = 8

# Row 40
$ws.Range("C40").Value = 4
$ws.Range("I40").Value = $null
$ws.Range("J40").Value = "new"
$ws.Range("K40").Value = "f"
$ws.Range("L40").Value = "stimuli/img_5il0t.png"
$ws.Range("M40").Value = 48.09523809523809
$ws.Range("N40").Value = 30.90476190476191
$ws.Range("O40").Value = 39.5
$ws.Range("Q40").Value = 2
$ws.Range("R40").Value = 2
$ws.Range("S40").Value = 2
$ws.Range("T40").Value = 2
$ws.Range("U40").Value = 2
$ws.Range("V40").Value = 2

# Row 41
$ws.Range("C41").Value = 4
$ws.Range("I41").Value = $null
$ws.Range("J41").Value = "new"
$ws.Range("K41").Value = "f"
$ws.Range("L41").Value = "stimuli/img_cogrz.png"
$ws.Range("M41").Value = 60.5
$ws.Range("N41").Value = 39.71428571428572
$ws.Range("O41").Value = 50.10714285714286
$ws.Range("P41").Value = 42
$ws.Range("Q41").Value = 3
$ws.Range("R41").Value = 3
$ws.Range("S41").Value = 3
$ws.Range("T41").Value = 3
$ws.Range("U41").Value = 3
$ws.Range("V41").Value = 3

# Row 42
$ws.Range("C42").Value = 4
$ws.Range("I42").Value = "target"
$ws.Range("J42").Value = "old"
$ws.Range("K42").Value = "j"
$ws.Range("L42").Value = "stimuli/img_fnu4h.png"
$ws.Range("M42").Value = 85.87179487179488
$ws.Range("N42").Value = 70.71794871794872
$ws.Range("O42").Value = 78.2948717948718
$ws.Range("P42").Value = 39
$ws.Range("Q42").Value = 9
$ws.Range("R42").Value = 9
$ws.Range("S42").Value = 9
$ws.Range("T42").Value = 9
$ws.Range("U42").Value = 9
$ws.Range("V42").Value = 9
